# Updates Price (D) and Volume(1h) (E) columns for the cryptos list refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.046.96"
$ws.Range("E2").Value = "  -3.56%  "

$ws.Range("D3").Value = "2.270.28"
$ws.Range("E3").Value = "  -4.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "490.82"
$ws.Range("E5").Value = "  -1.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.81"
$ws.Range("E6").Value = "  -2.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.525"
$ws.Range("E8").Value = "  -3.82%  "

$ws.Range("D9").Value = "2.270.35"
$ws.Range("E9").Value = "  -4.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0928"
$ws.Range("E10").Value = "  -4.25%  "

$ws.Range("E11").Value = "  -1.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.75"
$ws.Range("E12").Value = "  +2.59%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.316"
$ws.Range("E13").Value = "  -3.29%  "

$ws.Range("D14").Value = "2.673.87"
$ws.Range("E14").Value = "  -4.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.11"
$ws.Range("E15").Value = "  -1.20%  "

$ws.Range("D16").Value = "54.007.90"
$ws.Range("E16").Value = "  -3.48%  "

$ws.Range("E17").Value = "  -2.39%  "

$ws.Range("D18").Value = "2.314.73"
$ws.Range("E18").Value = "  -0.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.99"
$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.63"
$ws.Range("E20").Value = "  -3.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "302.70"
$ws.Range("E21").Value = "  -1.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.14"
$ws.Range("E22").Value = "  -1.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.92"
$ws.Range("E24").Value = "  -1.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.365"
$ws.Range("E26").Value = "  -0.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.144"
$ws.Range("E27").Value = "  -2.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.05"
$ws.Range("E28").Value = "  -2.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.09"
$ws.Range("E29").Value = "  -1.76%  "

$ws.Range("D30").Value = "0.0₃0695"
$ws.Range("E30").Value = "  -2.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.61"
$ws.Range("E31").Value = "  -1.27%  "

$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.77"
$ws.Range("E33").Value = "  +0.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  +0.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.07"
$ws.Range("E35").Value = "  -1.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.53"
$ws.Range("E36").Value = "  -0.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.834"
$ws.Range("E38").Value = "  +4.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.61"
$ws.Range("E39").Value = "  -4.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.77"
$ws.Range("E40").Value = "  -0.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.37"
$ws.Range("E41").Value = "  -2.41%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.366"
$ws.Range("E42").Value = "  -0.99%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.31"
$ws.Range("E43").Value = "  -0.82%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.58"
$ws.Range("E44").Value = "  -6.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.69"
$ws.Range("E45").Value = "  -1.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0880"
$ws.Range("E46").Value = "  -2.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.541"
$ws.Range("E47").Value = "  -4.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "238.20"
$ws.Range("E48").Value = "  -1.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0473"
$ws.Range("E49").Value = "  -1.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0204"
$ws.Range("E50").Value = "  -1.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.37"
$ws.Range("E51").Value = "  -2.73%  "

$ws.Range("D2:D51").Style = "Normal"
